$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has two picture "logos" duplicated across the header/footer
# parts. This change only renames the display name (wp:docPr / pic:cNvPr
# "name" attribute) of three of those inline pictures - the picture bytes,
# ids, alt-text/description and everything else stay the same:
#
#   footer (type "first",   Footers.Item(2))  Pearson logo  image2.png -> image1.png
#   footer (type "default", Footers.Item(1))  Pearson logo  image2.png -> image1.png
#   header (type "first",   Headers.Item(2))  BTEC logo     image1.jpg -> image2.jpg
#
# Renaming an InlineShape directly off a freshly-fetched
# Footers(...).Range.InlineShapes collection can leave a stale handle, so
# the shape is Select()-ed first and then addressed through the active
# Selection, which re-seats it reliably before the rename.

function Rename-Logo($shape, [string]$newName) {
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

$footerFirst = $sec.Footers.Item(2)
Rename-Logo $footerFirst.Range.InlineShapes.Item(1) "image1.png"

$footerDefault = $sec.Footers.Item(1)
Rename-Logo $footerDefault.Range.InlineShapes.Item(1) "image1.png"

$headerFirst = $sec.Headers.Item(2)
Rename-Logo $headerFirst.Range.InlineShapes.Item(1) "image2.jpg"
